$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-18 down to 11-19
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new data record
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "Vega Monumental Concepción"
$ws.Range("C10").Value = "Bíobío"
$ws.Range("D10").Value = 44539
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 100112022
$ws.Range("G10").Value = "Arveja Verde"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13400
$ws.Range("N10").Value = "$/saco 25 kilos"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 536
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
